# Generate Report for Handoff
# Mark the file "cac4ab63-da83-4c89-bb54-73862573db5f.md" as "Ready for handoff"
# after a fresh handoff xliff generation for the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the cac4ab63-... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F3").Value = "Ready for handoff"   # de-de status
$wsOverview.Range("G3").Value = "2016-08-21 10:13:31" # Latest HO Xliff Generate Date

# --- zh-cn sheet: row 3 is the cac4ab63-... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"       # Status
$wsZhCn.Range("E3").Value = "mt"                      # Priority
$wsZhCn.Range("H3").Value = "2016-08-21 10:13:27"     # Latest Handoff Datetime

# --- de-de sheet: row 3 is the cac4ab63-... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"       # Status
$wsDeDe.Range("E3").Value = "mt"                      # Priority
$wsDeDe.Range("H3").Value = "2016-08-21 10:13:31"     # Latest Handoff Datetime
